$d = $word.ActiveDocument
$d.Content.Find.Execute("MySQL, HTML, CSS, C#, PL/SQL, PL/pgSQL, Python, C++,React Native, Node.js, Spanish (Fluent)", $true, $false, $false, $false, $false, $true, 1, $false, "MySQL, HTML, CSS, C#, PL/SQL (read and debug), PL/pgSQL (read and debug), Python, C++, React Native, Node.js, Spanish (Fluent)", 2)
